# The VISITES table (starting at row 121) contained a "TYPVISIT" /
# "type de visite" field (row 126) that was removed from the data
# dictionary. Deleting the whole row shifts every row below it up by
# one, which is exactly what the target workbook shows (VOLS table and
# everything after it moves from rows 133-147 to rows 132-146, the
# sheet dimension shrinks from AMJ147 to AMJ146, and the now-unused
# "TYPVISIT" / "type de visite" shared strings disappear).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the TYPVISIT row from the VISITES table, shifting rows below it up.
$ws.Rows(126).EntireRow.Delete() | Out-Null

# Reflect the author's final selection/scroll position in the sheet view.
$ws.Range("J134").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 109
$excel.ActiveWindow.ScrollColumn = 1
